$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.983.46"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "1.677.38"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.24"
$ws.Range("E5").Value = "  +7.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3652"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.94"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3243"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.143"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07068"
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9981"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.093"
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.68"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "1.670.45"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.631"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001048"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06547"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9985"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.87"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.89"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.920"
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.86"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").Value = "24.955.47"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.444"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.399"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.03"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.75"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("D29").Value = "1.852.74"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.52"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.182"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.089"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.799"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08443"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.641"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.31"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.163"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06050"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02239"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.224"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2087"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.244"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9979"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5964"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  +9.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.852"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5730"
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.63"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.966"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07004"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.195"
$ws.Range("E51").Value = "  +3.87%  "
